# Update the "Exchange Conversion Date/Time" column (K) from 05/08/2025 to
# 06/08/2025 for every data row on the active sheet.
#
# These cells hold the date as literal text (not a real Excel date), so the
# new value must be written as text too. Assigning a date-looking string
# directly to .Value2 makes Excel auto-convert it into a date serial number,
# so we force text entry by prefixing with an apostrophe (same as typing
# '06/08/2025 into the cell) and then reset the cell style back to Normal so
# no stray "quote prefix" / text-number-format style sticks to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, "K").End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 30 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, "K")
    if ($cell.Value2 -eq "05/08/2025") {
        $cell.Value2 = "'06/08/2025"
        $cell.Style = "Normal"
    }
}
